$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell (untouched, default style) used as a formatting donor so that
# forcing text storage on numeric-looking strings does not leave a stray
# number-format style on the edited cell (matches the source file's styling).
$fmtDonor = $ws.Range("Z1")

$ws.Range("D2").Value = '60.203.62'
$ws.Range("E2").Value = '  +0.79%  '

$ws.Range("D3").Value = '2.414.85'
$ws.Range("E3").Value = '  -0.09%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '560.52'
$fmtDonor.Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("E5").Value = '  +1.63%  '

$ws.Range("E6").Value = '  -0.66%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("E8").Value = '  -0.20%  '

$ws.Range("E9").Value = '  +0.52%  '

$ws.Range("E10").Value = '  -0.92%  '

$ws.Range("E11").Value = '  +0.11%  '

$ws.Range("E12").Value = '  -0.95%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.86'
$fmtDonor.Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("E13").Value = '  -0.08%  '

$ws.Range("D14").Value = '2.846.08'
$ws.Range("E14").Value = '  -0.10%  '

$ws.Range("D15").Value = '60.132.38'
$ws.Range("E15").Value = '  +0.74%  '

$ws.Range("E16").Value = '  +0.71%  '

$ws.Range("D17").Value = '2.392.16'
$ws.Range("E17").Value = '  -0.88%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.53'
$fmtDonor.Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("E19").Value = '  +3.37%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '327.55'
$fmtDonor.Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").Value = '  -0.89%  '

$ws.Range("E21").Value = '  +1.82%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$fmtDonor.Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = '  +0.04%  '

$ws.Range("E23").Value = '  -1.63%  '

$ws.Range("E24").Value = '  +1.94%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.61'
$fmtDonor.Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("E25").Value = '  -1.29%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$fmtDonor.Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").Value = '  -0.22%  '

$ws.Range("E27").Value = '  +2.07%  '

$ws.Range("E28").Value = '  +2.33%  '

$ws.Range("D29").Value = '0.0₃0773'
$ws.Range("E29").Value = '  -0.13%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '170.43'
$fmtDonor.Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = '  +0.14%  '

$ws.Range("E31").Value = '  +0.17%  '

$ws.Range("E32").Value = '  +7.42%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.403'
$fmtDonor.Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("E33").Value = '  -2.00%  '

$ws.Range("E34").Value = '  -1.10%  '

$ws.Range("E35").Value = '  +4.20%  '

$ws.Range("E37").Value = '  -0.04%  '

$ws.Range("E38").Value = '  +0.54%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '324.66'
$fmtDonor.Copy()
$ws.Range("D39").PasteSpecial(-4122)
$ws.Range("E39").Value = '  +3.58%  '

$ws.Range("E40").Value = '  -0.15%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '38.58'
$fmtDonor.Copy()
$ws.Range("D41").PasteSpecial(-4122)
$ws.Range("E41").Value = '  -2.46%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '147.23'
$fmtDonor.Copy()
$ws.Range("D42").PasteSpecial(-4122)
$ws.Range("E42").Value = '  +5.97%  '

$ws.Range("E43").Value = '  -2.25%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0970'
$fmtDonor.Copy()
$ws.Range("D44").PasteSpecial(-4122)
$ws.Range("E44").Value = '  -0.04%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '19.79'
$fmtDonor.Copy()
$ws.Range("D45").PasteSpecial(-4122)
$ws.Range("E45").Value = '  +1.39%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0518'
$fmtDonor.Copy()
$ws.Range("D46").PasteSpecial(-4122)
$ws.Range("E46").Value = '  -0.03%  '

$ws.Range("E47").Value = '  -0.28%  '

$ws.Range("E48").Value = '  -1.06%  '

$ws.Range("E49").Value = '  -0.12%  '

$ws.Range("E50").Value = '  +0.07%  '

$ws.Range("E51").Value = '  -0.82%  '

$excel.CutCopyMode = $false
